$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update recalculated values in row 2
$ws.Range("M2").Value = 30.58864766666666
$ws.Range("N2").Value = 91.76594299999999
$ws.Range("O2").Value = 0.3925391465174898
$ws.Range("P2").Value = 0.3925391465174898
$ws.Range("Q2").Value = 12.916331775079
$ws.Range("R2").Value = 116.246985975711
$ws.Range("S2").Value = 0.3925391465174898
$ws.Range("T2").Value = 0.3925391465174898

# Update recalculated values in row 3
$ws.Range("O3").Value = 0.291183949679193
$ws.Range("P3").Value = 0.291183949679193
$ws.Range("S3").Value = 0.291183949679193
$ws.Range("T3").Value = 0.291183949679193

# Update recalculated values in row 4
$ws.Range("M4").Value = 24.64590566666666
$ws.Range("N4").Value = 73.93771699999999
$ws.Range("O4").Value = 0.3162769038033173
$ws.Range("P4").Value = 0.3162769038033172
$ws.Range("Q4").Value = 10.406955480901
$ws.Range("R4").Value = 93.662599328109
$ws.Range("S4").Value = 0.3162769038033173
$ws.Range("T4").Value = 0.3162769038033172

# Remove row 5 (Resolving-Mac entry) entirely
$ws.Rows.Item(5).Delete()

$wb.Save()
